$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.250.22'
$ws.Cells.Item(2, 5).Value = '  -0.57%  '

$ws.Cells.Item(3, 4).Value = '1.803.22'
$ws.Cells.Item(3, 5).Value = '  -0.81%  '

$ws.Cells.Item(4, 4).Value = '1.004'
$ws.Cells.Item(4, 5).Value = '  +0.24%  '

$ws.Cells.Item(5, 4).Value = '313.98'
$ws.Cells.Item(5, 5).Value = '  -0.52%  '

$ws.Cells.Item(6, 5).Value = '  +0.21%  '

$ws.Cells.Item(7, 4).Value = '0.5263'
$ws.Cells.Item(7, 5).Value = '  +2.83%  '

$ws.Cells.Item(8, 4).Value = '0.3820'
$ws.Cells.Item(8, 5).Value = '  -3.46%  '

$ws.Cells.Item(9, 4).Value = '0.08031'
$ws.Cells.Item(9, 5).Value = '  -1.58%  '

$ws.Cells.Item(10, 2).Value = 'Polygon'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(10, 4).Value = '1.101'
$ws.Cells.Item(10, 5).Value = '  -0.57%  '

$ws.Cells.Item(11, 2).Value = 'OKB'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(11, 4).Value = '41.31'
$ws.Cells.Item(11, 5).Value = '  -0.83%  '

$ws.Cells.Item(12, 4).Value = '6.316'
$ws.Cells.Item(12, 5).Value = '  +0.66%  '

$ws.Cells.Item(13, 4).Value = '1.004'
$ws.Cells.Item(13, 5).Value = '  +0.25%  '

$ws.Cells.Item(14, 4).Value = '20.59'
$ws.Cells.Item(14, 5).Value = '  -1.89%  '

$ws.Cells.Item(15, 4).Value = '1.807.84'
$ws.Cells.Item(15, 5).Value = '  -0.56%  '

$ws.Cells.Item(16, 4).Value = '7.319'
$ws.Cells.Item(16, 5).Value = '  -2.52%  '

$ws.Cells.Item(17, 4).Value = '92.23'
$ws.Cells.Item(17, 5).Value = '  -0.40%  '

$ws.Cells.Item(18, 4).Value = '0.00001095'
$ws.Cells.Item(18, 5).Value = '  -3.76%  '

$ws.Cells.Item(19, 4).Value = '0.06609'
$ws.Cells.Item(19, 5).Value = '  -0.37%  '

$ws.Cells.Item(20, 5).Value = '  +0.24%  '

$ws.Cells.Item(21, 4).Value = '17.35'
$ws.Cells.Item(21, 5).Value = '  -1.90%  '

$ws.Cells.Item(22, 4).Value = '5.969'
$ws.Cells.Item(22, 5).Value = '  -2.12%  '

$ws.Cells.Item(23, 4).Value = '28.306.58'
$ws.Cells.Item(23, 5).Value = '  -0.50%  '

$ws.Cells.Item(24, 4).Value = '11.13'
$ws.Cells.Item(24, 5).Value = '  -1.28%  '

$ws.Cells.Item(25, 4).Value = '2.277'
$ws.Cells.Item(25, 5).Value = '  +0.69%  '

$ws.Cells.Item(26, 4).Value = '160.56'
$ws.Cells.Item(26, 5).Value = '  +3.34%  '

$ws.Cells.Item(27, 4).Value = '20.46'
$ws.Cells.Item(27, 5).Value = '  -3.26%  '

$ws.Cells.Item(28, 4).Value = '2.011.55'
$ws.Cells.Item(28, 5).Value = '  -0.84%  '

$ws.Cells.Item(29, 4).Value = '2.357'
$ws.Cells.Item(29, 5).Value = '  -2.13%  '

$ws.Cells.Item(30, 4).Value = '123.31'
$ws.Cells.Item(30, 5).Value = '  -2.16%  '

$ws.Cells.Item(31, 4).Value = '0.1084'
$ws.Cells.Item(31, 5).Value = '  -1.60%  '

$ws.Cells.Item(32, 4).Value = '1.056'
$ws.Cells.Item(32, 5).Value = '  -4.51%  '

$ws.Cells.Item(33, 4).Value = '3.686'
$ws.Cells.Item(33, 5).Value = '  +0.90%  '

$ws.Cells.Item(34, 4).Value = '5.551'
$ws.Cells.Item(34, 5).Value = '  -3.81%  '

$ws.Cells.Item(35, 4).Value = '0.07227'
$ws.Cells.Item(35, 5).Value = '  +2.86%  '

$ws.Cells.Item(36, 4).Value = '12.36'
$ws.Cells.Item(36, 5).Value = '  +9.31%  '

$ws.Cells.Item(37, 4).Value = '0.02311'
$ws.Cells.Item(37, 5).Value = '  -0.72%  '

$ws.Cells.Item(38, 4).Value = '0.2147'
$ws.Cells.Item(38, 5).Value = '  -3.58%  '

$ws.Cells.Item(39, 4).Value = '5.120'
$ws.Cells.Item(39, 5).Value = '  -2.02%  '

$ws.Cells.Item(40, 4).Value = '8.610'
$ws.Cells.Item(40, 5).Value = '  -2.51%  '

$ws.Cells.Item(41, 4).Value = '0.6198'
$ws.Cells.Item(41, 5).Value = '  -1.27%  '

$ws.Cells.Item(42, 4).Value = '1.169'
$ws.Cells.Item(42, 5).Value = '  -0.47%  '

$ws.Cells.Item(43, 4).Value = '1.371'
$ws.Cells.Item(43, 5).Value = '  -2.10%  '

$ws.Cells.Item(44, 4).Value = '13.21'
$ws.Cells.Item(44, 5).Value = '  -2.18%  '

$ws.Cells.Item(45, 4).Value = '0.6014'
$ws.Cells.Item(45, 5).Value = '  +1.39%  '

$ws.Cells.Item(46, 5).Value = '  +0.71%  '

$ws.Cells.Item(47, 4).Value = '127.16'
$ws.Cells.Item(47, 5).Value = '  +1.90%  '

$ws.Cells.Item(48, 4).Value = '1.219'
$ws.Cells.Item(48, 5).Value = '  +2.71%  '

$ws.Cells.Item(49, 4).Value = '1.927'
$ws.Cells.Item(49, 5).Value = '  -2.50%  '

$ws.Cells.Item(50, 5).Value = '  -1.01%  '

$ws.Cells.Item(51, 4).Value = '73.05'
$ws.Cells.Item(51, 5).Value = '  -1.76%  '
